$d = $word.ActiveDocument

# Locate the paragraph that starts the "After cloning..." sentence.
$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "^After cloning") {
        $targetIdx = $i
        break
    }
}
if ($targetIdx -eq -1) {
    throw "Could not find the 'After cloning' paragraph"
}

# The four pieces of text that should end up as four separate runs
# inside the (now un-numbered) "After cloning" paragraph.
$piece1 = "After cloning, the selected columns in the target data table will have the same properties as the corresponding columns in the source data table with two except"
$piece2 = " that t"
$piece3 = "he column names are not changed. I.e., the target data table columns keep the names they have in the target data table" + [char]0x2019 + "s CSV file"
$piece4 = ", so the metadata will match the CSV file."

$p = $d.Paragraphs.Item($targetIdx)
$r = $p.Range
$r.Text = $piece1

# Build the remaining three pieces as separate trailing paragraphs, then
# splice the paragraph breaks back out so the runs stay distinct instead
# of being coalesced into a single run (which happens if text is merely
# appended into the same run via InsertAfter / Find-Replace).
$insAfter = $d.Range($d.Paragraphs.Item($targetIdx).Range.End - 1, $d.Paragraphs.Item($targetIdx).Range.End - 1)
$insAfter.InsertParagraphAfter()
$d.Paragraphs.Item($targetIdx + 1).Range.Text = $piece2

$insAfter2 = $d.Range($d.Paragraphs.Item($targetIdx + 1).Range.End - 1, $d.Paragraphs.Item($targetIdx + 1).Range.End - 1)
$insAfter2.InsertParagraphAfter()
$d.Paragraphs.Item($targetIdx + 2).Range.Text = $piece3

$insAfter3 = $d.Range($d.Paragraphs.Item($targetIdx + 2).Range.End - 1, $d.Paragraphs.Item($targetIdx + 2).Range.End - 1)
$insAfter3.InsertParagraphAfter()
$d.Paragraphs.Item($targetIdx + 3).Range.Text = $piece4

# Re-join: delete the 3 paragraph marks we just inserted so the 4 runs
# end up living inside one paragraph again.
for ($k = 0; $k -lt 3; $k++) {
    $mark = $d.Range($d.Paragraphs.Item($targetIdx).Range.End - 1, $d.Paragraphs.Item($targetIdx).Range.End)
    $mark.Delete()
}

Write-Host "After step 1: idx=$targetIdx text=[$($d.Paragraphs.Item($targetIdx).Range.Text)]"
Write-Host "Count: $($d.Paragraphs.Count)"

# Step 2: delete the now-empty paragraph right after the merged one
# (was the blank line that used to separate "exceptions:" from the
# bulleted list).
$d.Paragraphs.Item($targetIdx + 1).Range.Delete()

# Step 3: the old "column names" bullet paragraph's text has already
# been folded into the paragraph above, so just delete this whole
# paragraph (including its own paragraph mark).
$d.Paragraphs.Item($targetIdx + 1).Range.Delete()

# Step 4: the following paragraph used to be a blank ListParagraph
# bullet-list spacer (<w:pStyle w:val="ListParagraph"/><w:ind
# w:firstLine="0"/>). It becomes a plain paragraph (no style) with
# <w:ind w:left="0" w:firstLine="0"/>.
$spacer = $d.Paragraphs.Item($targetIdx + 1)
$spacer.Style = "Normal"
$spacer.LeftIndent = 0
$spacer.FirstLineIndent = 0

# Step 5: delete the "For Categorical columns..." bullet paragraph and
# the blank ListParagraph bullet spacer that follows it -- the whole
# second bullet point goes away (cloning a categorical column's codes
# is no longer restricted to those already present in the target).
$d.Paragraphs.Item($targetIdx + 2).Range.Delete()
$d.Paragraphs.Item($targetIdx + 2).Range.Delete()

Write-Host "Final count: $($d.Paragraphs.Count)"
